# "clean-up of input tables" - a cosmetic re-save of the scenario workbook
# (the commit only touches view/UI bookkeeping; every cell value, the
# shared-string table and the table definition are unchanged).
#
# Two things are adjusted through the Excel object model:
#   1. the visible sheet tab name: "updated" -> "Tabelle1"
#   2. which cell is selected/active when the sheet is (re)opened:
#      AN15 -> B10 (this also brings the scrolled-in top-left corner back
#      to column A, since B10 is already on-screen without scrolling).
#
# (A handful of other attributes in the diff - the absPath cache written
# by the Windows copy of Excel, the window's on-screen position/size, and
# sub-pixel bestFit column-width / font-metric caching values such as
# x14ac:dyDescent - are internal bookkeeping Excel stamps on save from its
# own environment; they have no corresponding read/write surface on the
# Application/Workbook/Worksheet object model, so they aren't set here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab.
$ws.Name = "Tabelle1"

# Make it the active sheet, then select B10 so that becomes the stored
# active cell / selection for the sheet view.
$ws.Activate() | Out-Null
$ws.Range("B10").Select() | Out-Null
